# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on row 2 of the
# zh-cn and de-de worksheets to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 15:32:14"
$wsZhCn.Range("H2").Value = "2016-03-24 15:32:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 15:32:18"
$wsDeDe.Range("H2").Value = "2016-03-24 15:32:54"
